$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "DATE"
$ws.Range("B1").Value = "ANNOUNCEMENT"

# --- Data row ---
$ws.Range("A2").Value = 45920
$ws.Range("A2").NumberFormat = "m/d/yyyy"
$ws.Range("B2").Value = "Song practice this week will be as follows:`n09/23/2025: 6:30PM to 8:30M`n09/25/2025: 6:00PM to 8:00PM"
$ws.Rows.Item(2).RowHeight = 60

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 10.140625
$ws.Columns.Item(2).ColumnWidth = 54.85546875

# --- Turn the range into an Excel Table ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:B2"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleLight2"

Write-Output "done"
